# Applies the "Add files via upload" edit:
#   - Slide 1, shape 1 (title "Cím 1"): reposition/resize, change text to
#     "Számítógépes vírusok", switch font to Century Schoolbook.
#   - Slide 1, shape 2 (subtitle "Alcím 2"): reposition/resize, bump font
#     size to 28pt, switch font to Century Schoolbook, drop the
#     lnSpcReduction on normAutofit.
#   - Slide 10, shape 1 ("Szövegdoboz 4"): reposition/resize, bump font
#     size to 36pt, switch font from Baskerville Old Face to Century
#     Schoolbook.
#
# EMU/point note: this host's Shape.Left/Top/Width/Height setters convert
# points -> EMU by truncation (floor), not rounding, so a point value
# computed as emu/12700 can land 1 EMU short. Using (emu + 0.5)/12700
# lands on the correct EMU after the floor.

$p = $ppt.ActivePresentation

# ---- Slide 1 : title placeholder ------------------------------------
$slide1 = $p.Slides.Item(1)

$title = $slide1.Shapes.Item(1)
$title.Left   = 186.00003937007875   # 2362200 EMU
$title.Top    = 193.00011811023623   # 2451101 EMU
$title.Width  = 573.0000393700788    # 7277100 EMU
$title.Height = 55.99996062992126    # 711199 EMU

$titleRange = $title.TextFrame.TextRange
$titleRange.Text = "Számítógépes vírusok"
$titleRange.Font.Name = "Century Schoolbook"

# ---- Slide 1 : subtitle placeholder ----------------------------------
$subtitle = $slide1.Shapes.Item(2)
$subtitle.Left   = 284.50003937007875  # 3613150 EMU
$subtitle.Top    = 277.00003937007875  # 3517900 EMU
$subtitle.Width  = 376.00003937007875  # 4775200 EMU
$subtitle.Height = 74.00003937007874   # 939800 EMU

# ppAutoSizeTextToFitShape -> plain <a:normAutofit/> (no lnSpcReduction)
$subtitle.TextFrame2.AutoSize = 2

$subtitleRange = $subtitle.TextFrame.TextRange
$subtitleRange.Font.Size = 28
$subtitleRange.Font.Name = "Century Schoolbook"

# ---- Slide 10 : "Köszönöm a figyelmet!" textbox ----------------------
$slide10 = $p.Slides.Item(10)
$thanks = $slide10.Shapes.Item(1)

$thanks.Left   = 271.00003937007875   # 3441700 EMU
$thanks.Top    = 229.00003937007875   # 2908300 EMU
$thanks.Width  = 391.00003937007875   # 4965700 EMU
$thanks.Height = 50.89224409448819    # 646331 EMU

$thanksRange = $thanks.TextFrame.TextRange
$thanksRange.Font.Size = 36
$thanksRange.Font.Name = "Century Schoolbook"
